$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected range (columns B-E, the coin name/link/price/volume
# columns) to Text format so numeric-looking strings (e.g. "1.005",
# "0.5123") are preserved exactly as literal text, matching the original
# inline-string cell contents instead of being auto-converted to numbers.
$affected = $ws.Range("B2:E51")
$affected.NumberFormat = "@"

$ws.Range("D2").Value = "26.784.47"
$ws.Range("E2").Value = "  -7.10%  "
$ws.Range("D3").Value = "1.696.08"
$ws.Range("E3").Value = "  -6.44%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "220.48"
$ws.Range("E5").Value = "  -5.18%  "
$ws.Range("D6").Value = "0.5123"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -6.57%  "
$ws.Range("D9").Value = "21.96"
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("D10").Value = "0.06190"
$ws.Range("E10").Value = "  -8.24%  "
$ws.Range("D11").Value = "0.07357"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "1.695.53"
$ws.Range("E12").Value = "  -6.44%  "
$ws.Range("D13").Value = "4.476"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "0.5811"
$ws.Range("E14").Value = "  -6.87%  "
$ws.Range("D15").Value = "1.926.97"
$ws.Range("E15").Value = "  -6.40%  "
$ws.Range("D16").Value = "0.000008191"
$ws.Range("E16").Value = "  -12.00%  "
$ws.Range("D17").Value = "65.28"
$ws.Range("E17").Value = "  -12.64%  "
$ws.Range("D18").Value = "26.818.49"
$ws.Range("E18").Value = "  -6.33%  "
$ws.Range("D19").Value = "5.024"
$ws.Range("E19").Value = "  -7.89%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "10.66"
$ws.Range("E21").Value = "  -6.32%  "
$ws.Range("D22").Value = "187.05"
$ws.Range("E22").Value = "  -10.47%  "
$ws.Range("D23").Value = "6.266"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "142.52"
$ws.Range("E25").Value = "  -7.69%  "
$ws.Range("D26").Value = "7.492"
$ws.Range("E26").Value = "  -4.04%  "
$ws.Range("D27").Value = "0.1147"
$ws.Range("E27").Value = "  -9.75%  "
$ws.Range("D28").Value = "15.21"
$ws.Range("E28").Value = "  -6.86%  "
$ws.Range("D29").Value = "1.334"
$ws.Range("E29").Value = "  -5.28%  "
$ws.Range("D30").Value = "0.05883"
$ws.Range("E30").Value = "  -6.50%  "
$ws.Range("D31").Value = "1.352"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("D32").Value = "3.465"
$ws.Range("E32").Value = "  -7.30%  "
$ws.Range("D33").Value = "3.443"
$ws.Range("E33").Value = "  -6.85%  "
$ws.Range("D34").Value = "1.651"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").Value = "0.9930"
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").Value = "2.415"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("D37").Value = "0.5971"
$ws.Range("D38").Value = "2.669"
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").Value = "1.091.73"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "0.01597"
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").Value = "5.843"
$ws.Range("E42").Value = "  -9.24%  "
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "97.27"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "1.842.18"
$ws.Range("E45").Value = "  -6.41%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "56.07"
$ws.Range("E46").Value = "  -7.43%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000103"
$ws.Range("E48").Value = "  -7.44%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.024"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05235"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4322"

# Remove the temporary Text-format styling so the saved cell styles match
# the original workbook (values remain text since they are already stored
# as strings at this point).
$affected.ClearFormats()

Write-Host "Applied cryptos list update."
